$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 15.83809466666666
$ws.Range("H2").Value = 47.514284
$ws.Range("I2").Value = 0.216761684821562
$ws.Range("J2").Value = 0.216761684821562
$ws.Range("M2").Value = 8.586181333333334
$ws.Range("N2").Value = 25.758544
$ws.Range("O2").Value = 0.2026609900758243
$ws.Range("P2").Value = 0.2026609900758243
$ws.Range("Q2").Value = 135.9887527824995
$ws.Range("R2").Value = 1223.898775042496
$ws.Range("S2").Value = 0.04392913765644153
$ws.Range("T2").Value = 0.04392913765644153
# Row 3
$ws.Range("G3").Value = 15.83809466666666
$ws.Range("H3").Value = 47.514284
$ws.Range("I3").Value = 0.216761684821562
$ws.Range("J3").Value = 0.216761684821562
$ws.Range("O3").Value = 0.5436514315916885
$ws.Range("P3").Value = 0.5436514315916886
$ws.Range("Q3").Value = 364.7987711049546
$ws.Range("R3").Value = 3283.188939944591
$ws.Range("S3").Value = 0.1178428002674685
$ws.Range("T3").Value = 0.1178428002674686
# Row 4
$ws.Range("G4").Value = 15.83809466666666
$ws.Range("H4").Value = 47.514284
$ws.Range("I4").Value = 0.216761684821562
$ws.Range("J4").Value = 0.216761684821562
$ws.Range("M4").Value = 10.74803566666667
$ws.Range("N4").Value = 32.244107
$ws.Range("O4").Value = 0.2536875783324871
$ws.Range("P4").Value = 0.2536875783324871
$ws.Range("Q4").Value = 170.2284063693764
$ws.Range("R4").Value = 1532.055657324388
$ws.Range("S4").Value = 0.05498974689765188
$ws.Range("T4").Value = 0.05498974689765188
# Row 5
$ws.Range("I5").Value = 0.4135749676850743
$ws.Range("J5").Value = 0.4135749676850743
$ws.Range("M5").Value = 8.586181333333334
$ws.Range("N5").Value = 25.758544
$ws.Range("O5").Value = 0.2026609900758243
$ws.Range("P5").Value = 0.2026609900758243
$ws.Range("Q5").Value = 259.4625709975165
$ws.Range("R5").Value = 2335.163138977648
$ws.Range("S5").Value = 0.08381551242163419
$ws.Range("T5").Value = 0.08381551242163419
# Row 6
$ws.Range("I6").Value = 0.4135749676850743
$ws.Range("J6").Value = 0.4135749676850743
$ws.Range("O6").Value = 0.5436514315916885
$ws.Range("P6").Value = 0.5436514315916886
$ws.Range("S6").Value = 0.2248406232524769
$ws.Range("T6").Value = 0.224840623252477
# Row 7
$ws.Range("I7").Value = 0.4135749676850743
$ws.Range("J7").Value = 0.4135749676850743
$ws.Range("M7").Value = 10.74803566666667
$ws.Range("N7").Value = 32.244107
$ws.Range("O7").Value = 0.2536875783324871
$ws.Range("P7").Value = 0.2536875783324871
$ws.Range("Q7").Value = 324.7908306361965
$ws.Range("R7").Value = 2923.117475725769
$ws.Range("S7").Value = 0.1049188320109631
$ws.Range("T7").Value = 0.1049188320109631
# Row 8
$ws.Range("G8").Value = 27.01013833333333
$ws.Range("H8").Value = 81.030415
$ws.Range("I8").Value = 0.3696633474933637
$ws.Range("J8").Value = 0.3696633474933637
$ws.Range("M8").Value = 8.586181333333334
$ws.Range("N8").Value = 25.758544
$ws.Range("O8").Value = 0.2026609900758243
$ws.Range("P8").Value = 0.2026609900758243
$ws.Range("Q8").Value = 231.9139455684178
$ws.Range("R8").Value = 2087.22551011576
$ws.Range("S8").Value = 0.07491633999774856
$ws.Range("T8").Value = 0.07491633999774858
# Row 9
$ws.Range("G9").Value = 27.01013833333333
$ws.Range("H9").Value = 81.030415
$ws.Range("I9").Value = 0.3696633474933637
$ws.Range("J9").Value = 0.3696633474933637
$ws.Range("O9").Value = 0.5436514315916885
$ws.Range("P9").Value = 0.5436514315916886
$ws.Range("Q9").Value = 622.1244081911133
$ws.Range("R9").Value = 5599.11967372002
$ws.Range("S9").Value = 0.200968008071743
$ws.Range("T9").Value = 0.200968008071743
# Row 10
$ws.Range("G10").Value = 27.01013833333333
$ws.Range("H10").Value = 81.030415
$ws.Range("I10").Value = 0.3696633474933637
$ws.Range("J10").Value = 0.3696633474933637
$ws.Range("M10").Value = 10.74803566666667
$ws.Range("N10").Value = 32.244107
$ws.Range("O10").Value = 0.2536875783324871
$ws.Range("P10").Value = 0.2536875783324871
$ws.Range("Q10").Value = 290.3059301682672
$ws.Range("R10").Value = 2612.753371514405
$ws.Range("S10").Value = 0.09377899942387209
$ws.Range("T10").Value = 0.09377899942387212

$wb.Save()
